$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 111473773
$ws.Range("B3").Value = 89405
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 1202
$ws.Range("F3").Value = "Ullticka"
$ws.Range("G3").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H3").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q3").Value = 704016.0051346947
$ws.Range("R3").Value = 6572801.994589122
$ws.Range("AO3").Value = "granlåga"

# Row 4
$ws.Range("A4").Value = 111473779
$ws.Range("B4").Value = 89425
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 5442
$ws.Range("F4").Value = "Tallticka"
$ws.Range("G4").Value = "Porodaedalea pini"
$ws.Range("H4").Value = "(Brot.) Murrill"
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value = 704193.4830821306
$ws.Range("R4").Value = 6572948.378178579
$ws.Range("AO4").Value = "gammeltall"

# Row 5
$ws.Range("A5").Value = 111473777
$ws.Range("B5").Value = 89425
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5442
$ws.Range("F5").Value = "Tallticka"
$ws.Range("G5").Value = "Porodaedalea pini"
$ws.Range("H5").Value = "(Brot.) Murrill"
$ws.Range("Q5").Value = 704301.1177162804
$ws.Range("R5").Value = 6573209.392206083
$ws.Range("AO5").Value = "gammeltall"

# Row 6
$ws.Range("A6").Value = 111473792
$ws.Range("B6").Value = 5113
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 100526
$ws.Range("F6").Value = "Bronshjon"
$ws.Range("G6").Value = "Callidium coriaceum"
$ws.Range("H6").Value = "Paykull, 1800"
$ws.Range("M6").Value = "äldre gnagspår"
$ws.Range("Q6").Value = 703965.55072247
$ws.Range("R6").Value = 6572785.445717536
$ws.Range("AO6").Value = "torrgran"

# Row 7
$ws.Range("A7").Value = 111473784
$ws.Range("B7").Value = 73634
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 6426
$ws.Range("F7").Value = "Kattfotslav"
$ws.Range("G7").Value = "Felipes leucopellaeus"
$ws.Range("H7").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q7").Value = 704135.470341172
$ws.Range("R7").Value = 6572843.267234835
$ws.Range("AO7").Value = "äldre gran"

# Row 9
$ws.Range("A9").Value = 111473774
$ws.Range("Q9").Value = 703999.5190368021
$ws.Range("R9").Value = 6572850.823973293

# Row 10
$ws.Range("A10").Value = 111473782
$ws.Range("B10").Value = 89183
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 3215
$ws.Range("F10").Value = "Rödgul trumpetsvamp"
$ws.Range("G10").Value = "Craterellus lutescens"
$ws.Range("H10").Value = "(Fr.) Fr."
$ws.Range("Q10").Value = 704171.5165585374
$ws.Range("R10").Value = 6572850.843097115
$ws.Range("AO10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111473791
$ws.Range("B11").Value = 93289
$ws.Range("E11").Value = 2170
$ws.Range("F11").Value = "Flagellkvastmossa"
$ws.Range("G11").Value = "Dicranum flagellare"
$ws.Range("H11").Value = "Hedw."
$ws.Range("Q11").Value = 704004.9502936595
$ws.Range("R11").Value = 6572835.740028554
$ws.Range("AO11").Value = "låga av tall"

# Row 12
$ws.Range("A12").Value = 111473775
$ws.Range("Q12").Value = 703969.3444121893
$ws.Range("R12").Value = 6572791.287347207

# Row 13
$ws.Range("A13").Value = 111473776
$ws.Range("B13").Value = 89405
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = "Ullticka"
$ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q13").Value = 703970.8884549731
$ws.Range("R13").Value = 6572810.333898042
$ws.Range("AO13").Value = "granlåga"
